# Update "想去人数" (want-to-go count) figures in column F for two
# sheets that list the same events: "展览" and "全部类型".
#
#   Sheet "展览"    : F2 123->125, F3 432->438, F5 8->9
#   Sheet "全部类型" : F2 123->125, F4 432->438, F6 8->9

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 125
$wsExhibit.Range("F3").Value = 438
$wsExhibit.Range("F5").Value = 9

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 125
$wsAll.Range("F4").Value = 438
$wsAll.Range("F6").Value = 9
